$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unicode helpers (avoid relying on literal multi-byte characters in the script)
$omega = [char]0x03A9
$pm    = [char]0x00B1

# --- Fix the resistor description typo (0? -> 0Ω) used by rows 4 and 6 ---
$resistor0ohm = "RESISTOR: SMD 0603 0" + $omega + " 1A [RC0603FR-070RL] [SMD]"
$ws.Range("C4").Value = $resistor0ohm
$ws.Range("C6").Value = $resistor0ohm

# --- Row 5: DNI group now only covers R2, R5, R6, R7 ---
$ws.Range("B5").Value = "R2, R5, R6, R7"
$ws.Range("E5").Value = "-"

# --- Row 7 becomes a populated resistor row (R8, R9, R10, R11) ---
$ws.Range("A7").Value = "OEPS020013"
$ws.Range("B7").Value = "R8, R9, R10, R11"
$ws.Range("C7").Value = "RESISTOR: SMD 0402 10k" + $omega + " " + $pm + "1% 0.100W [ERJ-2RKF1002X] [SMD]"
$ws.Range("D7").Value = "DNP"
$ws.Range("E7").Value = "ERJ-2RKF1002X"

# --- Former rows 7 & 8 (connectors) shift down to rows 8 & 9 ---
$ws.Range("A8").Value = "OEPS070053"
$ws.Range("B8").Value = "U$4, U$5, U$6"
$ws.Range("C8").Value = "CONNECTOR: CONN RCPT 4POS 0.079 GOLD PCB [37204-62A3-004PL] [TH]"
$ws.Range("E8").Value = "37204-62A3-004PL"

$ws.Range("A9").Value = "OEPS070054"
$ws.Range("B9").Value = "U$8"
$ws.Range("C9").Value = "CONNECTOR: CONN RCPT 12POS 0.079 GOLD PCB [37212-62M3-003PL] [TH]"
$ws.Range("E9").Value = "37212-62M3-003PL"

# --- New note + cable row ---
$ws.Range("A12").Value = "POST-PRODUCTION"
$ws.Range("A12").Font.Bold = $true

$ws.Range("A13").Value = "OEPS090005"
$ws.Range("C13").Value = "CABLE: FFC / FPC 26WAY 0.5MM 152MM [MP-FFCA05261522A]"
$ws.Range("E13").Value = "MP-FFCA05261522A"

# --- Column widths (C widened, new F column width) ---
$ws.Columns("C").ColumnWidth = 71.9
$ws.Columns("F").ColumnWidth = 24.8

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection as saved in the file ---
$ws.Range("C17").Select() | Out-Null
